$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06446533333333333
$ws.Range("H2").Value = 0.193396
$ws.Range("I2").Value = 0.02693738696927793
$ws.Range("J2").Value = 0.02693738696927793
$ws.Range("M2").Value = 7.757543333333333
$ws.Range("N2").Value = 23.27263
$ws.Range("O2").Value = 0.4040769763164727
$ws.Range("P2").Value = 0.4040769763164727
$ws.Range("Q2").Value = 0.5000926168311112
$ws.Range("R2").Value = 4.50083355148
$ws.Range("S2").Value = 0.01088477787641258
$ws.Range("T2").Value = 0.01088477787641258

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06446533333333333
$ws.Range("H3").Value = 0.193396
$ws.Range("I3").Value = 0.02693738696927793
$ws.Range("J3").Value = 0.02693738696927793
$ws.Range("O3").Value = 0.01627055103446774
$ws.Range("P3").Value = 0.01627055103446774
$ws.Range("Q3").Value = 0.02013671384666667
$ws.Range("R3").Value = 0.18123042462
$ws.Range("S3").Value = 0.0004382861294188429
$ws.Range("T3").Value = 0.0004382861294188429

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06446533333333333
$ws.Range("H4").Value = 0.193396
$ws.Range("I4").Value = 0.02693738696927793
$ws.Range("J4").Value = 0.02693738696927793
$ws.Range("M4").Value = 11.12827366666667
$ws.Range("N4").Value = 33.384821
$ws.Range("O4").Value = 0.5796524726490594
$ws.Range("P4").Value = 0.5796524726490595
$ws.Range("Q4").Value = 0.7173878713462222
$ws.Range("R4").Value = 6.456490842116001
$ws.Range("S4").Value = 0.0156143229634465
$ws.Range("T4").Value = 0.01561432296344651

# Row 5
$ws.Range("I5").Value = 0.7704314695358874
$ws.Range("J5").Value = 0.7704314695358874
$ws.Range("M5").Value = 7.757543333333333
$ws.Range("N5").Value = 23.27263
$ws.Range("O5").Value = 0.4040769763164727
$ws.Range("P5").Value = 0.4040769763164727
$ws.Range("Q5").Value = 14.30306102550556
$ws.Range("R5").Value = 128.72754922955
$ws.Range("S5").Value = 0.3113136186691181
$ws.Range("T5").Value = 0.3113136186691181

# Row 6
$ws.Range("I6").Value = 0.7704314695358874
$ws.Range("J6").Value = 0.7704314695358874
$ws.Range("O6").Value = 0.01627055103446774
$ws.Range("P6").Value = 0.01627055103446774
$ws.Range("S6").Value = 0.01253534454364364
$ws.Range("T6").Value = 0.01253534454364364

# Row 7
$ws.Range("I7").Value = 0.7704314695358874
$ws.Range("J7").Value = 0.7704314695358874
$ws.Range("M7").Value = 11.12827366666667
$ws.Range("N7").Value = 33.384821
$ws.Range("O7").Value = 0.5796524726490594
$ws.Range("P7").Value = 0.5796524726490595
$ws.Range("Q7").Value = 20.51788440277611
$ws.Range("R7").Value = 184.660959624985
$ws.Range("S7").Value = 0.4465825063231256
$ws.Range("T7").Value = 0.4465825063231257

# Row 8
$ws.Range("G8").Value = 0.4849276666666666
$ws.Range("H8").Value = 1.454783
$ws.Range("I8").Value = 0.2026311434948347
$ws.Range("J8").Value = 0.2026311434948347
$ws.Range("M8").Value = 7.757543333333333
$ws.Range("N8").Value = 23.27263
$ws.Range("O8").Value = 0.4040769763164727
$ws.Range("P8").Value = 0.4040769763164727
$ws.Range("Q8").Value = 3.761847387698889
$ws.Range("R8").Value = 33.85662648929
$ws.Range("S8").Value = 0.0818785797709421
$ws.Range("T8").Value = 0.0818785797709421

# Row 9
$ws.Range("G9").Value = 0.4849276666666666
$ws.Range("H9").Value = 1.454783
$ws.Range("I9").Value = 0.2026311434948347
$ws.Range("J9").Value = 0.2026311434948347
$ws.Range("O9").Value = 0.01627055103446774
$ws.Range("P9").Value = 0.01627055103446774
$ws.Range("Q9").Value = 0.1514744305983333
$ws.Range("R9").Value = 1.363269875385
$ws.Range("S9").Value = 0.003296920361405264
$ws.Range("T9").Value = 0.003296920361405264

# Row 10
$ws.Range("G10").Value = 0.4849276666666666
$ws.Range("H10").Value = 1.454783
$ws.Range("I10").Value = 0.2026311434948347
$ws.Range("J10").Value = 0.2026311434948347
$ws.Range("M10").Value = 11.12827366666667
$ws.Range("N10").Value = 33.384821
$ws.Range("O10").Value = 0.5796524726490594
$ws.Range("P10").Value = 0.5796524726490595
$ws.Range("Q10").Value = 5.396407783204777
$ws.Range("R10").Value = 48.567670048843
$ws.Range("S10").Value = 0.1174556433624873
$ws.Range("T10").Value = 0.1174556433624873

